$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, shifting existing rows 62-141 down to 63-142.
$ws.Rows("62:62").Insert()

# Populate the newly inserted row 62 with the new record's data.
$ws.Range("A62").Value = 2
$ws.Range("B62").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44902
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 100112024
$ws.Range("G62").Value = "Choclo"
$ws.Range("H62").Value = "Dulce o Americano"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 11000
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = 11500
$ws.Range("N62").Value = "$/malla 70 unidades"
$ws.Range("O62").Value = "Provincia de Limarí"
$ws.Range("P62").Value = 164
$ws.Range("Q62").Value = 70
$ws.Range("R62").Value = "Hortaliza"
